# P-197 Ajout des fichier .c et .h dans le repertoire CCS
# dsk.h/dsk.c
# wave.h/wave.c
#
# Updates the MFCC pipeline performance worksheet:
#  - fill in the previously-empty F8/F9 "% d'un cycle" cells
#  - move the stray G12/G16/G21 cells back into column F
#  - add a total in F23 (sum of F8:F22)
#  - add a new "Pipeline MFCC 13 coefficient (256 donne)" row (row 25)
#  - refresh the sheet view (scroll position / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the two cells in the first shared-formula block that were left blank ---
$ws.Range("F8").Formula = "=E8/10"
$ws.Range("F9").Formula = "=E9/10"

# --- The G12/G16/G21 cells belong in column F ("% d'un cycle du MFCC") ---
$ws.Range("F12").NumberFormat = $ws.Range("G12").NumberFormat
$ws.Range("F12").Formula = "=E12/10"
$ws.Range("G12").Clear()

$ws.Range("F16").NumberFormat = $ws.Range("G16").NumberFormat
$ws.Range("F16").Formula = "=E16/10"
$ws.Range("G16").Clear()

$ws.Range("F21").NumberFormat = $ws.Range("G21").NumberFormat
$ws.Range("F21").Formula = "=E21/10"
$ws.Range("G21").Clear()

# --- Total line right under the table ---
$ws.Range("F23").Formula = "=SUM(F8:F22)"

# --- New row: Pipeline MFCC 13 coefficient (256 donne) ---
$ws.Range("C25").Value = "Pipeline MFCC 13 coefficient (256 donné)"
$ws.Range("D25").Value = 78783
$ws.Range("E25").Formula = "=D25/225000000*1000"
$ws.Range("F25").Formula = "=E25/10"

# --- Refresh the view: scroll back to the top, select D24 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("D24").Select()

$wb.Application.Calculate()
